$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bracket")

$ws.Range('D2').Value = 'Okapi'
$ws.Range('N2').Value = 'Golden Eagle'
$ws.Range('E4').Value = 'Okapi'
$ws.Range('M4').Value = 'Golden Eagle'
$ws.Range('D6').Value = 'Striped Polecat'
$ws.Range('N6').Value = 'Veined Octopus'
$ws.Range('F8').Value = 'Okapi'
$ws.Range('L8').Value = 'Golden Eagle'
$ws.Range('D10').Value = 'Side-striped jackal'
$ws.Range('N10').Value = 'Palaeocastor'
$ws.Range('E12').Value = 'Striped dolphin'
$ws.Range('M12').Value = 'Lungfish'
$ws.Range('D14').Value = 'Striped dolphin'
$ws.Range('N14').Value = 'Lungfish'
$ws.Range('G16').Value = 'Okapi'
$ws.Range('K16').Value = 'Golden Eagle'
$ws.Range('D18').Value = 'Wildcat'
$ws.Range('N18').Value = 'Goanna'
$ws.Range('E20').Value = 'Striped hyena'
$ws.Range('M20').Value = 'Homo habilis'
$ws.Range('D22').Value = 'Striped hyena'
$ws.Range('N22').Value = 'Homo habilis'
$ws.Range('F24').Value = 'Kudu'
$ws.Range('L24').Value = 'Cathedral Termite'
$ws.Range('D26').Value = 'Striped Rabbit'
$ws.Range('N26').Value = 'Montezuma Oropendola'
$ws.Range('E28').Value = 'Kudu'
$ws.Range('M28').Value = 'Cathedral Termite'
$ws.Range('D30').Value = 'Kudu'
$ws.Range('N30').Value = 'Cathedral Termite'
$ws.Range('H32').Value = 'Okapi'
$ws.Range('I32').Value = 'Golden Eagle'
$ws.Range('J32').Value = 'Golden Eagle'
$ws.Range('D34').Value = 'Sea Otter'
$ws.Range('N34').Value = 'Emperor Penguin'
$ws.Range('E36').Value = 'Sea Otter'
$ws.Range('M36').Value = 'Emperor Penguin'
$ws.Range('D38').Value = 'Southern Ningaui'
$ws.Range('N38').Value = 'Owl Monkey'
$ws.Range('F40').Value = 'Sea Otter'
$ws.Range('L40').Value = 'Emperor Penguin'
$ws.Range('D42').Value = 'Sibree Dwarf Lemur'
$ws.Range('N42').Value = 'Pacific Spiny Lumpsucker'
$ws.Range('E44').Value = 'Mara'
$ws.Range('M44').Value = 'Siamang'
$ws.Range('D46').Value = 'Mara'
$ws.Range('N46').Value = 'Siamang'
$ws.Range('G48').Value = 'Sea Otter'
$ws.Range('K48').Value = 'Emperor Penguin'
$ws.Range('D50').Value = 'Itjaritjari'
$ws.Range('N50').Value = 'Bat-Eared Fox'
$ws.Range('E52').Value = 'Dik Dik'
$ws.Range('M52').Value = 'Wolverine'
$ws.Range('D54').Value = 'Dik Dik'
$ws.Range('N54').Value = 'Wolverine'
$ws.Range('F56').Value = 'Rock Hyrax '
$ws.Range('L56').Value = 'Greater Rhea'
$ws.Range('D58').Value = 'Bulldog Bat'
$ws.Range('N58').Value = 'Greater Flamingo'
$ws.Range('E60').Value = 'Rock Hyrax '
$ws.Range('M60').Value = 'Greater Rhea'
$ws.Range('D62').Value = 'Rock Hyrax '
$ws.Range('N62').Value = 'Greater Rhea'
